$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.31951409500135
$ws.Range("C2").Value = 8.421901884717794
$ws.Range("E2").Value = 15.41338964439735
$ws.Range("F2").Value = 38.70716184524979
$ws.Range("G2").Value = 3.664789623265651
$ws.Range("I2").Value = 25.64012644928269
$ws.Range("J2").Value = 8.128711292047736
$ws.Range("K2").Value = 9.915423429260676
$ws.Range("L2").Value = 12.58501652232382
$ws.Range("N2").Value = 19.82368920817019
$ws.Range("O2").Value = 25.46630810466185

$ws.Range("B3").Value = 13.0957736436648
$ws.Range("C3").Value = 8.413835677159037
$ws.Range("E3").Value = 15.4076678102597
$ws.Range("F3").Value = 38.72701607147268
$ws.Range("G3").Value = 3.666400352119401
$ws.Range("I3").Value = 25.72824997208459
$ws.Range("J3").Value = 8.122333459163174
$ws.Range("K3").Value = 9.76491553122958
$ws.Range("L3").Value = 12.56470262702102
$ws.Range("N3").Value = 19.8806308285916
$ws.Range("O3").Value = 25.54830361004298

$ws.Range("B4").Value = 12.95864625596074
$ws.Range("C4").Value = 8.408916072292596
$ws.Range("E4").Value = 15.40645347949633
$ws.Range("F4").Value = 38.74747177045469
$ws.Range("G4").Value = 3.667442258338199
$ws.Range("I4").Value = 25.78646307715764
$ws.Range("J4").Value = 8.118403366025538
$ws.Range("K4").Value = 9.672783890907898
$ws.Range("L4").Value = 12.5540417396292
$ws.Range("N4").Value = 19.917258321135
$ws.Range("O4").Value = 25.60313661622176

$ws.Range("B5").Value = 12.90290002437177
$ws.Range("C5").Value = 8.406919950105705
$ws.Range("E5").Value = 15.40653878926022
$ws.Range("F5").Value = 38.75788667002743
$ws.Range("G5").Value = 3.667880189158277
$ws.Range("I5").Value = 25.81121757868561
$ws.Range("J5").Value = 8.116798428200104
$ws.Range("K5").Value = 9.635354697415071
$ws.Range("L5").Value = 12.55015637104398
$ws.Range("N5").Value = 19.93260422276991
$ws.Range("O5").Value = 25.62660904319126

$ws.Range("B6").Value = 12.89365361575102
$ws.Range("C6").Value = 8.40658902824347
$ws.Range("E6").Value = 15.40658805126372
$ws.Range("F6").Value = 38.75974162751809
$ws.Range("G6").Value = 3.667953714465225
$ws.Range("I6").Value = 25.815390386051
$ws.Range("J6").Value = 8.11653172881755
$ws.Range("K6").Value = 9.629147872842998
$ws.Range("L6").Value = 12.54953902675058
$ws.Range("N6").Value = 19.93517779426266
$ws.Range("O6").Value = 25.63057469610058

$ws.Range("B7").Value = 12.95789380420907
$ws.Range("C7").Value = 8.408889116390197
$ws.Range("E7").Value = 15.40645227839023
$ws.Range("F7").Value = 38.7476038112455
$ws.Range("G7").Value = 3.667448110336645
$ws.Range("I7").Value = 25.78679274575551
$ws.Range("J7").Value = 8.118381734819526
$ws.Range("K7").Value = 9.672278582231618
$ws.Range("L7").Value = 12.5539874773054
$ws.Range("N7").Value = 19.91746357958107
$ws.Range("O7").Value = 25.60344860942048

$ws.Range("B8").Value = 13.24235904181496
$ws.Range("C8").Value = 8.419113992731566
$ws.Range("E8").Value = 15.41094093622436
$ws.Range("F8").Value = 38.71229283983186
$ws.Range("G8").Value = 3.665334044992904
$ws.Range("I8").Value = 25.6696594536228
$ws.Range("J8").Value = 8.126515224699526
$ws.Range("K8").Value = 9.863496682945705
$ws.Range("L8").Value = 12.57763835051813
$ws.Range("N8").Value = 19.84297780882661
$ws.Range("O8").Value = 25.49364827379184

$ws.Range("B9").Value = 13.79909424219039
$ws.Range("C9").Value = 8.439413267308066
$ws.Range("E9").Value = 15.43788404713001
$ws.Range("F9").Value = 38.70856267365831
$ws.Range("G9").Value = 3.661606360198013
$ws.Range("I9").Value = 25.47253290876476
$ws.Range("J9").Value = 8.142347317428099
$ws.Range("K9").Value = 10.2387643449703
$ws.Range("L9").Value = 12.63823401461629
$ws.Range("N9").Value = 19.71006788354083
$ws.Range("O9").Value = 25.31397410853214

$ws.Range("B10").Value = 14.20343384289995
$ws.Range("C10").Value = 8.454460614525081
$ws.Range("E10").Value = 15.46859203891056
$ws.Range("F10").Value = 38.74562192037857
$ws.Range("G10").Value = 3.659119866406749
$ws.Range("I10").Value = 25.34756147733002
$ws.Range("J10").Value = 8.153900180865161
$ws.Range("K10").Value = 10.51210371856222
$ws.Range("L10").Value = 12.69118965344785
$ws.Range("N10").Value = 19.6203606794436
$ws.Range("O10").Value = 25.20374809216451

$ws.Range("B11").Value = 14.38554879511267
$ws.Range("C11").Value = 8.461331213302726
$ws.Range("E11").Value = 15.48489536053045
$ws.Range("F11").Value = 38.77107344697674
$ws.Range("G11").Value = 3.658042917834075
$ws.Range("I11").Value = 25.29502114118871
$ws.Range("J11").Value = 8.159137168913835
$ws.Range("K11").Value = 10.63541956149579
$ws.Range("L11").Value = 12.71705751833669
$ws.Range("N11").Value = 19.58125823183878
$ws.Range("O11").Value = 25.15834336756936

$ws.Range("B12").Value = 14.45418505890274
$ws.Range("C12").Value = 8.463936239030108
$ws.Range("E12").Value = 15.49140113178826
$ws.Range("F12").Value = 38.78194095253365
$ws.Range("G12").Value = 3.65764285354876
$ws.Range("I12").Value = 25.27574538854339
$ws.Range("J12").Value = 8.16111746442791
$ws.Range("K12").Value = 10.68192686230067
$ws.Range("L12").Value = 12.72710358293857
$ws.Range("N12").Value = 19.56669521016185
$ws.Range("O12").Value = 25.14183191360069

$ws.Range("B13").Value = 14.43941854214322
$ws.Range("C13").Value = 8.463375060469716
$ws.Range("E13").Value = 15.48998528894409
$ws.Range("F13").Value = 38.77954586383481
$ws.Range("G13").Value = 3.657728670327463
$ws.Range("I13").Value = 25.27986917991762
$ws.Range("J13").Value = 8.160691101531599
$ws.Range("K13").Value = 10.67191977164119
$ws.Range("L13").Value = 12.72492892954626
$ws.Range("N13").Value = 19.5698207722654
$ws.Range("O13").Value = 25.14535758734656

$ws.Range("B14").Value = 14.3912024502787
$ws.Range("C14").Value = 8.461545466840102
$ws.Range("E14").Value = 15.48542396074496
$ws.Range("F14").Value = 38.77194291956309
$ws.Range("G14").Value = 3.658009849137889
$ws.Range("I14").Value = 25.29342288278051
$ws.Range("J14").Value = 8.159300146957941
$ws.Range("K14").Value = 10.63924978205073
$ws.Range("L14").Value = 12.71787902524136
$ws.Range("N14").Value = 19.58005523468547
$ws.Range("O14").Value = 25.15697127946485

$ws.Range("B15").Value = 14.36162428125464
$ws.Range("C15").Value = 8.460425203444515
$ws.Range("E15").Value = 15.48267314852214
$ws.Range("F15").Value = 38.76744582330448
$ws.Range("G15").Value = 3.65818308783003
$ws.Range("I15").Value = 25.30180568838265
$ws.Range("J15").Value = 8.158447772020804
$ws.Range("K15").Value = 10.61921254439708
$ws.Range("L15").Value = 12.7135932139189
$ws.Range("N15").Value = 19.58635591248111
$ws.Range("O15").Value = 25.16417389156205

$ws.Range("B16").Value = 14.19148994517163
$ws.Range("C16").Value = 8.454012084674002
$ws.Range("E16").Value = 15.46757324598889
$ws.Range("F16").Value = 38.74413101940242
$ws.Range("G16").Value = 3.659191334412903
$ws.Range("I16").Value = 25.35108186900369
$ws.Range("J16").Value = 8.153557544238522
$ws.Range("K16").Value = 10.50402043527674
$ws.Range("L16").Value = 12.6895344669282
$ws.Range("N16").Value = 19.62295035015127
$ws.Range("O16").Value = 25.20681082291755

$ws.Range("B17").Value = 14.08660452189859
$ws.Range("C17").Value = 8.450084098057687
$ws.Range("E17").Value = 15.45890529450717
$ws.Range("F17").Value = 38.73202501894671
$ws.Range("G17").Value = 3.659823708963458
$ws.Range("I17").Value = 25.3824153255226
$ws.Range("J17").Value = 8.150552689129643
$ws.Range("K17").Value = 10.43306023732648
$ws.Range("L17").Value = 12.67522698358858
$ws.Range("N17").Value = 19.64583598843642
$ws.Range("O17").Value = 25.23418139176282

$ws.Range("B18").Value = 14.02610947949793
$ws.Range("C18").Value = 8.447827246115983
$ws.Range("E18").Value = 15.4541396453265
$ws.Range("F18").Value = 38.72587122118495
$ws.Range("G18").Value = 3.66019253498904
$ws.Range("I18").Value = 25.40084313534886
$ws.Range("J18").Value = 8.148822648843243
$ws.Range("K18").Value = 10.39215149584521
$ws.Range("L18").Value = 12.66716538821087
$ws.Range("N18").Value = 19.65915982415721
$ws.Range("O18").Value = 25.25037008726189

$ws.Range("B19").Value = 14.00560011587873
$ws.Range("C19").Value = 8.44706354250981
$ws.Range("E19").Value = 15.45256395618767
$ws.Range("F19").Value = 38.72392681035118
$ws.Range("G19").Value = 3.660318290357187
$ws.Range("I19").Value = 25.40715213152742
$ws.Range("J19").Value = 8.148236595803066
$ws.Range("K19").Value = 10.37828560612768
$ws.Range("L19").Value = 12.66446482628575
$ws.Range("N19").Value = 19.66369866566161
$ws.Range("O19").Value = 25.25592784669434

$ws.Range("B20").Value = 14.09778759196813
$ws.Range("C20").Value = 8.450501991050539
$ws.Range("E20").Value = 15.45980527626272
$ws.Range("F20").Value = 38.7332300096428
$ws.Range("G20").Value = 3.659755863955777
$ws.Range("I20").Value = 25.37903784115232
$ws.Range("J20").Value = 8.15087274114509
$ws.Range("K20").Value = 10.4406241417427
$ws.Range("L20").Value = 12.67673272124602
$ws.Range("N20").Value = 19.64338315716749
$ws.Range("O20").Value = 25.23122159446389

$ws.Range("B21").Value = 14.40537405384062
$ws.Range("C21").Value = 8.462082776935508
$ws.Range("E21").Value = 15.48675475127215
$ws.Range("F21").Value = 38.77414277237904
$ws.Range("G21").Value = 3.657927049989173
$ws.Range("I21").Value = 25.28942499834247
$ws.Range("J21").Value = 8.159708782534207
$ws.Range("K21").Value = 10.64885122533039
$ws.Range("L21").Value = 12.71994299835869
$ws.Range("N21").Value = 19.57704250518993
$ws.Range("O21").Value = 25.15354153006481

$ws.Range("B22").Value = 14.60446734704884
$ws.Range("C22").Value = 8.469670421696913
$ws.Range("E22").Value = 15.50630143605782
$ws.Range("F22").Value = 38.80804490951603
$ws.Range("G22").Value = 3.656776986540283
$ws.Range("I22").Value = 25.2344726180466
$ws.Range("J22").Value = 8.165467015955254
$ws.Range("K22").Value = 10.78381550165625
$ws.Range("L22").Value = 12.74964071423921
$ws.Range("N22").Value = 19.53510810678748
$ws.Range("O22").Value = 25.10675072195219

$ws.Range("B23").Value = 14.49840493454887
$ws.Range("C23").Value = 8.465619148903091
$ws.Range("E23").Value = 15.49569332238624
$ws.Range("F23").Value = 38.78929754567675
$ws.Range("G23").Value = 3.6573866758656
$ws.Range("I23").Value = 25.26347085297903
$ws.Range("J23").Value = 8.162395321333719
$ws.Range("K23").Value = 10.71189881190094
$ws.Range("L23").Value = 12.73365893252299
$ws.Range("N23").Value = 19.55735942446994
$ws.Range("O23").Value = 25.13135956521496

$ws.Range("B24").Value = 14.09273232844489
$ws.Range("C24").Value = 8.450313057246426
$ws.Range("E24").Value = 15.45939771622733
$ws.Range("F24").Value = 38.73268272195804
$ws.Range("G24").Value = 3.659786520279064
$ws.Range("I24").Value = 25.38056351293505
$ws.Range("J24").Value = 8.150728053429676
$ws.Range("K24").Value = 10.43720484861635
$ws.Range("L24").Value = 12.67605146612194
$ws.Range("N24").Value = 19.64449156321647
$ws.Range("O24").Value = 25.23255830785793

$ws.Range("B25").Value = 13.64902343010496
$ws.Range("C25").Value = 8.433897026669229
$ws.Range("E25").Value = 15.42866721642046
$ws.Range("F25").Value = 38.70256848981511
$ws.Range("G25").Value = 3.662570317001947
$ws.Range("I25").Value = 25.52237380789224
$ws.Range("J25").Value = 8.138077647342479
$ws.Range("K25").Value = 10.13748016234295
$ws.Range("L25").Value = 12.69118965344785
$ws.Range("N25").Value = 19.7446232404371
$ws.Range("O25").Value = 25.3587594109528

